$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Italian lesson description strings for column C (rows 65-112),
# plus column E gets a copy of column B's ("en") text for each row.
$rows = @(
    @{ Row = 65; Text = "sommare due numeri a una cifra a 18" },
    @{ Row = 66; Text = "sommare un numero ad una cifra ed un numero a due cifre fino a 19" },
    @{ Row = 67; Text = "sottrarre un numero ad una cifra a zero" },
    @{ Row = 68; Text = "sottrarre un numero a due cifre a zero" },
    @{ Row = 69; Text = "sommare e sottrarre numeri ad una cifra fino a 18" },
    @{ Row = 70; Text = "sommare e sottrarre numeri ad una ed a due cifre a 20 incluso lo zero" },
    @{ Row = 71; Text = "sommare numeri a due cifre e unità" },
    @{ Row = 72; Text = "sommare numeri a due cifre e decine" },
    @{ Row = 73; Text = "sommare due numeri a due cifre" },
    @{ Row = 74; Text = "sommare tre numeri ad una cifra" },
    @{ Row = 75; Text = "sottrarre un numero a due cifre e unità" },
    @{ Row = 76; Text = "sottrarre un numero a due cifre e decine" },
    @{ Row = 77; Text = "sottrarre due numeri a due cifre a zero" },
    @{ Row = 78; Text = "sottrarre tre numeri ad una cifra a zero" },
    @{ Row = 79; Text = "sommare e sottrarre tre numeri a due cifre" },
    @{ Row = 80; Text = "moltiplicare la tabellina del 2 fino a 12 volte" },
    @{ Row = 81; Text = "moltiplicare la tabellina del 10 fino a 12 volte" },
    @{ Row = 82; Text = "moltiplicare la tabellina del 5 fino a 12 volte" },
    @{ Row = 83; Text = "moltiplicare 2, 5, 10 con numeri a due cifre" },
    @{ Row = 84; Text = "sommare un numero a tre cifre e unità" },
    @{ Row = 85; Text = "sommare un numero a tre cifre e decine" },
    @{ Row = 86; Text = "sommare un numero a tre cifre e centinaia" },
    @{ Row = 87; Text = "sottrarre un numero a tre cifre e unità" },
    @{ Row = 88; Text = "sottrarre un numero a tre cifre e decine" },
    @{ Row = 89; Text = "sottrarre un numero a tre cifre e centinaia" },
    @{ Row = 90; Text = "sommare e sottrarre un numero a tre cifre e unità" },
    @{ Row = 91; Text = "sommare e sottrarre un numero a tre cifre e decine" },
    @{ Row = 92; Text = "sommare e sottrarre un numero a tre cifre e centinaia" },
    @{ Row = 93; Text = "moltiplicare la tabellina del 3 fino a 12 volte" },
    @{ Row = 94; Text = "moltiplicare la tabellina del 4 fino a 12 volte" },
    @{ Row = 95; Text = "dividere per 10 e 100" },
    @{ Row = 96; Text = "moltiplicare la tabellina dell'8 fino a 12 volte" },
    @{ Row = 97; Text = "sommare e sottrarre due numeri a tre cifre" },
    @{ Row = 98; Text = "sommare numeri a quattro cifre" },
    @{ Row = 99; Text = "sottrarre numeri a quattro cifre" },
    @{ Row = 100; Text = "sommare e sottrarre numeri a quattro cifre" },
    @{ Row = 101; Text = "moltiplicare la tabellina del 12 fino a 12 volte" },
    @{ Row = 102; Text = "moltiplicare numeri a due cifre per un numero ad una cifra" },
    @{ Row = 103; Text = "moltiplicare numeri a tre cifre per un numero ad una cifra" },
    @{ Row = 104; Text = "dividere un numero a due cifre per un numero ad una cifra" },
    @{ Row = 105; Text = "dividere un numero a tre cifre per un numero ad una cifra" },
    @{ Row = 106; Text = "moltiplicare numeri fino a quattro cifre per un numero ad una cifra" },
    @{ Row = 107; Text = "riconoscere i quadrati dei numeri ed i simboli matematici" },
    @{ Row = 108; Text = "riconoscere i cubi dei numeri ed i simboli matematici" },
    @{ Row = 109; Text = "dividere numeri fino a quattro cifre per un numero ad una cifra" },
    @{ Row = 110; Text = "risolvere un'operazione (+,-,x, ÷) fino a numeri di quattro cifre" },
    @{ Row = 111; Text = "risolvere due operazioni fino a numeri di quattro cifre" },
    @{ Row = 112; Text = "risolvere tre operazioni fino a numeri di quattro cifre" }
)

foreach ($item in $rows) {
    $row = $item.Row

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value2 = $item.Text
    $cCell.HorizontalAlignment = 1

    $eCell = $ws.Cells.Item($row, 5)
    $bCell = $ws.Cells.Item($row, 2)
    $eCell.Value2 = $bCell.Value2
    $eCell.Style = "Bad"
    $eCell.Borders.LineStyle = 1
    $eCell.Borders.Weight = -4138
    $eCell.Borders.Color = 13421772
    $eCell.HorizontalAlignment = -4131
    $eCell.WrapText = $true
}

# Row-height tweaks for rows whose content now wraps onto an extra line.
$ws.Rows.Item(70).RowHeight = 29.5
$ws.Rows.Item(92).RowHeight = 29.5
$ws.Rows.Item(106).RowHeight = 26.5
$ws.Rows.Item(107).RowHeight = 26.5
$ws.Rows.Item(109).RowHeight = 29.5

# Update the sheet view scroll position / active selection.
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("E57").Select()
